{"js": "// Replace each \"NNN\u00d7N=\" expression in the document's table cells with its\n// updated counterpart, per the commit's regenerated numbers.\nconst replacements = [\n  [\"570\u00d74=\", \"439\u00d73=\"],\n  [\"966\u00d77=\", \"661\u00d72=\"],\n  [\"450\u00d77=\", \"883\u00d77=\"],\n  [\"435\u00d74=\", \"950\u00d72=\"],\n  [\"269\u00d78=\", \"453\u00d76=\"],\n  [\"286\u00d74=\", \"665\u00d73=\"],\n  [\"387\u00d77=\", \"150\u00d75=\"],\n  [\"803\u00d79=\", \"716\u00d76=\"],\n  [\"982\u00d79=\", \"597\u00d79=\"],\n  [\"908\u00d72=\", \"823\u00d74=\"],\n  [\"452\u00d74=\", \"378\u00d77=\"],\n  [\"366\u00d75=\", \"408\u00d79=\"],\n  [\"909\u00d78=\", \"251\u00d77=\"],\n  [\"710\u00d78=\", \"937\u00d77=\"],\n  [\"921\u00d77=\", \"237\u00d75=\"],\n  [\"120\u00d79=\", \"613\u00d77=\"],\n  [\"640\u00d78=\", \"459\u00d78=\"],\n  [\"647\u00d74=\", \"732\u00d72=\"],\n  [\"829\u00d75=\", \"427\u00d73=\"],\n  [\"265\u00d75=\", \"833\u00d77=\"],\n  [\"700\u00d77=\", \"850\u00d73=\"],\n  [\"245\u00d72=\", \"223\u00d76=\"],\n  [\"281\u00d79=\", \"400\u00d77=\"],\n  [\"418\u00d77=\", \"382\u00d76=\"],\n  [\"350\u00d72=\", \"141\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"NNNxN=\" multiplication expression in the document's\n# table cells with its updated counterpart, per the commit's regenerated\n# numbers. Uses Word's Find/Replace (wdReplaceAll) scoped per-pair so that\n# each distinct \"before\" string is swapped for its matching \"after\" string.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"570\u00d74=\", \"439\u00d73=\"),\n    @(\"966\u00d77=\", \"661\u00d72=\"),\n    @(\"450\u00d77=\", \"883\u00d77=\"),\n    @(\"435\u00d74=\", \"950\u00d72=\"),\n    @(\"269\u00d78=\", \"453\u00d76=\"),\n    @(\"286\u00d74=\", \"665\u00d73=\"),\n    @(\"387\u00d77=\", \"150\u00d75=\"),\n    @(\"803\u00d79=\", \"716\u00d76=\"),\n    @(\"982\u00d79=\", \"597\u00d79=\"),\n    @(\"908\u00d72=\", \"823\u00d74=\"),\n    @(\"452\u00d74=\", \"378\u00d77=\"),\n    @(\"366\u00d75=\", \"408\u00d79=\"),\n    @(\"909\u00d78=\", \"251\u00d77=\"),\n    @(\"710\u00d78=\", \"937\u00d77=\"),\n    @(\"921\u00d77=\", \"237\u00d75=\"),\n    @(\"120\u00d79=\", \"613\u00d77=\"),\n    @(\"640\u00d78=\", \"459\u00d78=\"),\n    @(\"647\u00d74=\", \"732\u00d72=\"),\n    @(\"829\u00d75=\", \"427\u00d73=\"),\n    @(\"265\u00d75=\", \"833\u00d77=\"),\n    @(\"700\u00d77=\", \"850\u00d73=\"),\n    @(\"245\u00d72=\", \"223\u00d76=\"),\n    @(\"281\u00d79=\", \"400\u00d77=\"),\n    @(\"418\u00d77=\", \"382\u00d76=\"),\n    @(\"350\u00d72=\", \"141\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $null = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
